$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R1").Value = 0.84783393725096001
$ws.Range("AK1").Value = 0.85164221402122675
$ws.Range("AL1").Value = 0.74021726911777597
$ws.Range("AO1").Value = 0.66244753931254485
$ws.Range("BA2").Value = 0.71541485904666713
$ws.Range("BC2").Value = 0.60551002463644166
$ws.Range("P3").Value = 0.64371342225794703
$ws.Range("BA3").Value = 0.92735325295060289
$ws.Range("C4").Value = 0.6874721310360179
$ws.Range("E4").Value = 0.80560854869224052
$ws.Range("U4").Value = 0.76644337936639761
$ws.Range("Z5").Value = 0.85953001256127948
$ws.Range("BG5").Value = 0.78238296823810694
$ws.Range("AL6").Value = 0.88940613293482351
$ws.Range("Z7").Value = 0.852041042658926
$ws.Range("F8").Value = 0.92477905362752377
$ws.Range("AB8").Value = 0.86061346397804805
$ws.Range("BG8").Value = 0.80479141959398715
$ws.Range("J9").Value = 0.81456143039258
$ws.Range("AU10").Value = 0.96621892580383162
$ws.Range("I11").Value = 0.91951348797212762
$ws.Range("M11").Value = 0.98459877246573391
$ws.Range("Q12").Value = 0.68479822942764068
$ws.Range("BA13").Value = 0.72343862020649419
$ws.Range("U14").Value = 0.88614145937524524
$ws.Range("N15").Value = 0.99722154637918414
$ws.Range("O17").Value = 0.91796047049006435
$ws.Range("P18").Value = 0.73648036673825956
$ws.Range("AN18").Value = 0.70001267668346012
$ws.Range("BA19").Value = 0.97679869928827046
$ws.Range("BL19").Value = 0.59460241159717286
$ws.Range("N20").Value = 0.58396731919409572
$ws.Range("AS20").Value = 0.99804545865566885
$ws.Range("BK20").Value = 0.85181030088720755
$ws.Range("AM22").Value = 0.79109007997128644
$ws.Range("AP22").Value = 0.81593270686320341
$ws.Range("M23").Value = 0.93273267806474303
$ws.Range("Y24").Value = 0.96486129745642502
$ws.Range("AA24").Value = 0.73923458038901246
$ws.Range("L25").Value = 0.78848858632606222
$ws.Range("W25").Value = 0.83997343161025861
$ws.Range("AJ25").Value = 0.95491772702992594
$ws.Range("BJ25").Value = 0.8531221435851819
$ws.Range("BP25").Value = 0.65794910941596474
$ws.Range("AL26").Value = 0.8163208784076661
$ws.Range("AC27").Value = 0.66537065374190973
$ws.Range("BE27").Value = 0.88601125498558542
$ws.Range("K28").Value = 0.5920336855277728
$ws.Range("S28").Value = 0.61377405653323303
$ws.Range("BB28").Value = 0.97073574143022379
$ws.Range("BO29").Value = 0.7632923112768788
$ws.Range("AE30").Value = 0.59135805717161616
$ws.Range("AQ30").Value = 0.95106644976585919
$ws.Range("BF30").Value = 0.725642636477416
$ws.Range("AN31").Value = 0.95590361309752736
$ws.Range("H32").Value = 0.88207520326175226
$ws.Range("BK32").Value = 0.91558251535201851
$ws.Range("A33").Value = 0.9931843240158853
$ws.Range("C33").Value = 0.81323188332200924
$ws.Range("I34").Value = 0.60515086012786068
$ws.Range("Q34").Value = 0.76506593883684015
$ws.Range("AE34").Value = 0.88936895427063845
$ws.Range("BN34").Value = 0.64682852046821071
$ws.Range("BO34").Value = 0.91180163838120487
$ws.Range("AG35").Value = 0.88824576319212656
$ws.Range("BF35").Value = 0.97357710149402765
$ws.Range("BP35").Value = 0.82090729969313836
$ws.Range("W37").Value = 0.58099888715473835
$ws.Range("AB37").Value = 0.95558722071509661
$ws.Range("AE37").Value = 0.97510533385956188
$ws.Range("AI37").Value = 0.72226529778061677
$ws.Range("AY37").Value = 0.89062645446558708
$ws.Range("BC37").Value = 0.71983912598316913
$ws.Range("AT38").Value = 0.93763551120995525
$ws.Range("AN39").Value = 0.89273946442813767
$ws.Range("G40").Value = 0.84629484388658072
$ws.Range("BA40").Value = 0.94742033117314817
$ws.Range("Q41").Value = 0.80433296188499015
$ws.Range("AF42").Value = 0.83719980521663051
$ws.Range("AO42").Value = 0.84471200622366371
$ws.Range("BD42").Value = 0.78820986408494276
$ws.Range("G43").Value = 0.68340753189130554
$ws.Range("AS43").Value = 0.67494028025040365
$ws.Range("BP43").Value = 0.86392799756917082
$ws.Range("AS44").Value = 0.57792272894193353
$ws.Range("BO44").Value = 0.79421078071864926
$ws.Range("BO45").Value = 0.6914149246154202
$ws.Range("AD46").Value = 0.9022192013864383
$ws.Range("AW47").Value = 0.8560396150135271
$ws.Range("Q48").Value = 0.95465305618302154
$ws.Range("T49").Value = 0.91786075944063494
$ws.Range("AZ49").Value = 0.89764701694972326
$ws.Range("F50").Value = 0.92143643589492585
$ws.Range("AW51").Value = 0.93482304923351622
$ws.Range("E52").Value = 0.58939316076429771
$ws.Range("O52").Value = 0.91173294738318567
$ws.Range("AA52").Value = 0.95434677774033894
$ws.Range("AM53").Value = 0.99198768863052256
$ws.Range("AX53").Value = 0.95016194461531134
$ws.Range("L54").Value = 0.97172902936404926
$ws.Range("X54").Value = 0.81892659073428387
$ws.Range("O55").Value = 0.79864562045518273
$ws.Range("N56").Value = 0.92437514140331489
$ws.Range("AR56").Value = 0.6643954369269971
$ws.Range("AG57").Value = 0.99742506287551036
$ws.Range("BO57").Value = 0.78494542425325187
$ws.Range("G58").Value = 0.84849963591106348
$ws.Range("T58").Value = 0.9125651747263408
$ws.Range("AJ58").Value = 0.91844655045440371
$ws.Range("AV58").Value = 0.6824627114193067
$ws.Range("BH59").Value = 0.61784591902469288
$ws.Range("BI59").Value = 0.65865309446175813
$ws.Range("B60").Value = 0.9414071055760026
$ws.Range("E60").Value = 0.96126297118391901
$ws.Range("AC60").Value = 0.9443603388774513
$ws.Range("U61").Value = 0.77005666305835052
$ws.Range("BJ61").Value = 0.83614815673690523
$ws.Range("BN61").Value = 0.75582848271421788
$ws.Range("F62").Value = 0.83338715882048819
$ws.Range("AT62").Value = 0.97479615501673189
$ws.Range("S63").Value = 0.96531910436593615
$ws.Range("U63").Value = 0.95853715552404417
$ws.Range("AL63").Value = 0.92487762288282971
$ws.Range("BA63").Value = 0.80328793084327321
$ws.Range("J64").Value = 0.76268805005728779
$ws.Range("X64").Value = 0.84732843713804651
$ws.Range("AY64").Value = 0.85857350603891769
$ws.Range("AZ64").Value = 0.9166212629072511
$ws.Range("AU65").Value = 0.89230409677537093
$ws.Range("BN65").Value = 0.95827204142244815
$ws.Range("H66").Value = 0.89827688647477699
$ws.Range("AW67").Value = 0.9986001756477092
$ws.Range("BM67").Value = 0.73790458178794971
$ws.Range("J68").Value = 0.66523041718317244
$ws.Range("O68").Value = 0.66924201719848475
$ws.Range("AY68").Value = 0.743508836714879
